$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.970.85"
Set-TextValue "E2" "  -0.85%  "
Set-TextValue "D3" "1.639.66"
Set-TextValue "E3" "  -0.26%  "
Set-TextValue "E4" "  +0.49%  "
Set-TextValue "D5" "215.84"
Set-TextValue "E5" "  -0.54%  "
Set-TextValue "E6" "  +0.23%  "
Set-TextValue "E7" "  +0.35%  "
Set-TextValue "D8" "0.257"
Set-TextValue "E8" "  -0.55%  "
Set-TextValue "E9" "  +0.09%  "
Set-TextValue "D10" "19.55"
Set-TextValue "E10" "  -1.19%  "
Set-TextValue "D11" "0.0793"
Set-TextValue "E11" "  +0.36%  "
Set-TextValue "E12" "  -0.09%  "
Set-TextValue "D13" "1.863.72"
Set-TextValue "E13" "  -0.46%  "
Set-TextValue "D14" "1.622.97"
Set-TextValue "E14" "  -1.51%  "
Set-TextValue "D15" "0.544"
Set-TextValue "E15" "  -0.20%  "
Set-TextValue "D16" "0.0₃0763"
Set-TextValue "E16" "  -0.21%  "
Set-TextValue "D17" "62.91"
Set-TextValue "E17" "  -0.50%  "
Set-TextValue "D18" "25.924.34"
Set-TextValue "E18" "  -1.04%  "
Set-TextValue "E19" "  +0.44%  "
Set-TextValue "D20" "192.75"
Set-TextValue "E20" "  -1.06%  "
Set-TextValue "E21" "  -1.74%  "
Set-TextValue "D22" "9.93"
Set-TextValue "E22" "  -1.29%  "
Set-TextValue "E23" "  -0.30%  "
Set-TextValue "E24" "  +0.96%  "
Set-TextValue "B25" "Monero"
Set-TextValue "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" "144.53"
Set-TextValue "E25" "  +1.36%  "
Set-TextValue "B26" "Stellar"
Set-TextValue "C26" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D26" "0.131"
Set-TextValue "E26" "  +4.98%  "
Set-TextValue "E27" "  +0.40%  "
Set-TextValue "D28" "6.91"
Set-TextValue "E28" "  -0.30%  "
Set-TextValue "E29" "  -0.51%  "
Set-TextValue "E30" "  -0.12%  "
Set-TextValue "E31" "  -0.45%  "
Set-TextValue "E32" "  -1.91%  "
Set-TextValue "E33" "  -0.11%  "
Set-TextValue "D34" "1.54"
Set-TextValue "E34" "  -3.60%  "
Set-TextValue "E35" "  +2.02%  "
Set-TextValue "D36" "0.901"
Set-TextValue "E36" "  -1.18%  "
Set-TextValue "D37" "1.133.17"
Set-TextValue "E37" "  -0.01%  "
Set-TextValue "D38" "0.542"
Set-TextValue "E38" "  -1.82%  "
Set-TextValue "D39" "2.46"
Set-TextValue "E39" "  -1.43%  "
Set-TextValue "E40" "  -0.22%  "
Set-TextValue "E41" "  +0.27%  "
Set-TextValue "D42" "99.59"
Set-TextValue "E42" "  -0.89%  "
Set-TextValue "D43" "0.794"
Set-TextValue "E43" "  -0.41%  "
Set-TextValue "D44" "1.773.12"
Set-TextValue "E44" "  -0.48%  "
Set-TextValue "D45" "0.0₆0114"
Set-TextValue "E45" "  +2.34%  "
Set-TextValue "D46" "56.61"
Set-TextValue "E46" "  -0.47%  "
Set-TextValue "D47" "0.0531"
Set-TextValue "E47" "  +2.78%  "
Set-TextValue "E48" "  -0.95%  "
Set-TextValue "D49" "7.71"
Set-TextValue "E49" "  +0.15%  "
Set-TextValue "E50" "  -0.51%  "
Set-TextValue "E51" "  -0.90%  "
